# ---------------------------------------------------------------------------
# Applies the "Leitura dos arquivos csv" commit:
#  - Adds two new lookup tables (ses_cias, ses_ramos) to the "fields" sheet,
#    plus fills in the previously-"undefined" metadata for the ses_seguros
#    table that was already present (row 2), including its new table_fk
#    reference (column I) to ses_cias.
#  - Extends the existing conditional formatting / data validation so it
#    keeps covering the newly added rows (and a bit of head-room below them).
#  - Leaves the "fields" sheet as the active tab (it was "validations"
#    before).
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$wsV = $wb.Worksheets.Item("validations")
$wsF = $wb.Worksheets.Item("fields")

# ---------------------------------------------------------------------------
# 1. Cell values
# ---------------------------------------------------------------------------

# Row 2 - ses_seguros (already existed, metadata was "undefined"; now filled
# in, and a table_fk to ses_cias is added in column I)
$wsF.Range("A2").Value = "ses_seguros.csv"
$wsF.Range("B2").Value = "ses_seguros"
$wsF.Range("C2").Value = "coenti"
$wsF.Range("D2").Value = "number"
$wsF.Range("E2").Value = "integer"
$wsF.Range("F2").Value = "no"
$wsF.Range("G2").Value = "yes"
$wsF.Range("H2").Value = "yes"
$wsF.Range("I2").Value = "ses_cias"
$wsF.Range("N2").Value = "yes"

# Row 3 - ses_cias (new)
$wsF.Range("A3").Value = "ses_cias.csv"
$wsF.Range("B3").Value = "ses_cias"
$wsF.Range("C3").Value = "coenti"
$wsF.Range("D3").Value = "number"
$wsF.Range("E3").Value = "integer"
$wsF.Range("F3").Value = "no"
$wsF.Range("G3").Value = "yes"
$wsF.Range("H3").Value = "no"
$wsF.Range("N3").Value = "yes"

# Row 4 - ses_ramos (new)
$wsF.Range("A4").Value = "ses_ramos.csv"
$wsF.Range("B4").Value = "ses_ramos"
$wsF.Range("C4").Value = "coramo"
$wsF.Range("D4").Value = "number"
$wsF.Range("E4").Value = "integer"
$wsF.Range("F4").Value = "no"
$wsF.Range("G4").Value = "yes"
$wsF.Range("H4").Value = "no"
$wsF.Range("N4").Value = "yes"

# ---------------------------------------------------------------------------
# 2. Column widths (A & B now hold file/table names and need to be widened)
# ---------------------------------------------------------------------------
$wsF.Columns.Item(1).ColumnWidth = 13.65
$wsF.Columns.Item(2).ColumnWidth = 10.3

# ---------------------------------------------------------------------------
# 3. Conditional formatting - extend existing rules to cover the new rows,
#    and add matching rules for rows 3 and 4 (new rows get their own rule,
#    mirroring what Excel does when formatting is copied down).
# ---------------------------------------------------------------------------
$fcDH = $wsF.Range("D2:H2").FormatConditions.Item(1)
$fcN  = $wsF.Range("N2").FormatConditions.Item(1)
$fillColor = $fcDH.Interior.Color()

$fcDH.ModifyAppliesToRange($wsF.Range("D2:I2"))
$fcN.ModifyAppliesToRange($wsF.Range("N2:N23"))

$fcRow3 = $wsF.Range("D3:I3").FormatConditions.Add(1, 3, '="undefined"')
$fcRow3.Interior.Color = $fillColor

$fcRow4 = $wsF.Range("D4:I23").FormatConditions.Add(1, 3, '="undefined"')
$fcRow4.Interior.Color = $fillColor

# ---------------------------------------------------------------------------
# 4. Data validation - widen the list-validation ranges so newly added rows
#    (and a bit of room below them) keep the dropdowns.
# ---------------------------------------------------------------------------
$wsF.Range("D2").Validation.Delete()
$wsF.Range("D2:D23").Validation.Add(3, 1, 1, '"undefined,number,text,data"')

$wsF.Range("E2").Validation.Delete()
$wsF.Range("E2:E23").Validation.Add(3, 1, 1, '"undefined,integer,decimal"')

$wsF.Range("N2").Validation.Delete()
$wsF.Range("F2:H2").Validation.Delete()
$wsF.Range("N2:N23").Validation.Add(3, 1, 1, '"undefined,yes,no"')
$wsF.Range("F2:H23").Validation.Add(3, 1, 1, '"undefined,yes,no"')

# ---------------------------------------------------------------------------
# 5. Active sheet - "fields" becomes the selected/active tab instead of
#    "validations".
# ---------------------------------------------------------------------------
$wsF.Activate()
$wsF.Range("A1").Select()
